# "I just deleted row 49" (Biofourmis) from the vendor list on Sheet1.
#
# Deleting the row shifts the remaining cell values/styles up correctly,
# but this host's Hyperlinks collection does not re-target itself when a
# row is removed (the hyperlink refs + underlying relationships are left
# stale). So the hyperlinks on column C are rebuilt explicitly: remember
# the cell style used for linked cells, clear the existing hyperlinks,
# delete the row, then re-create a hyperlink for each (now shifted) URL
# that already sits in the cell's own text, restoring the original style
# afterwards (Hyperlinks.Add mutates the cell style as a side effect).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# The style applied to the linked cells in column C (index 2 in styles.xml).
$linkStyle = $ws.Range("C2").Style

# Drop the existing hyperlinks on column C so we can cleanly rebuild them
# once the row below has shifted into place.
$ws.Range("C2:C" + $lastRow).Hyperlinks.Delete()

# Delete the Biofourmis row; everything below moves up one row.
$ws.Rows.Item(49).Delete()

$newLastRow = $lastRow - 1

# Re-create the hyperlinks for column C (rows 2..newLastRow). Each cell's
# own text already holds the target URL, so reuse it verbatim (setting
# .Address after Add avoids the host's "bare domain" URL auto-slash), and
# restore the cell style Add() stomped on.
for ($i = 2; $i -le $newLastRow; $i++) {
    $cell = $ws.Cells.Item($i, 3)
    $url = $cell.Value()
    $h = $ws.Hyperlinks.Add($cell, $url)
    $h.Address = $url
    $cell.Style = $linkStyle
}
